$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the full content of row 10 and row 11 (two occurrence
# records got re-ordered). Read the current ("before") values of the
# cells that differ between the two rows, then write them back swapped.
# Only cells whose value actually changes are touched, so that cells
# which are identical between the two rows (and therefore untouched by
# the diff) are not re-written and risk unwanted auto-conversion (e.g.
# Excel turning a literal date-like text such as "2026-01-24" into a
# real date when it is written back through .Value).

$cols = @("A","B","E","F","G","H","Q","R","AC")

$row10 = @{}
$row11 = @{}
foreach ($col in $cols) {
    $row10[$col] = $ws.Range($col + "10").Value()
    $row11[$col] = $ws.Range($col + "11").Value()
}

foreach ($col in $cols) {
    $ws.Range($col + "10").Value = $row11[$col]
    $ws.Range($col + "11").Value = $row10[$col]
}

# Columns M, AH, AJ, AK, AM, AO held data only on row 10 before the
# edit; after the edit that data belongs to row 11 and row 10's cells
# become empty.
$moveCols = @("M","AH","AJ","AK","AM","AO")
foreach ($col in $moveCols) {
    $val = $ws.Range($col + "10").Value()
    $ws.Range($col + "11").Value = $val
    $ws.Range($col + "10").ClearContents()
}
